$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws1.Range("A1").Value = "test"
